$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / non-numeric-looking values: direct assignment keeps Text type ---
$ws.Range("D2").Value = "60.909.47"
$ws.Range("D3").Value = "3.387.24"
$ws.Range("E3").Value = "  -1.44%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E5").Value = "  -0.80%  "
$ws.Range("E6").Value = "  -2.30%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.388.37"
$ws.Range("E8").Value = "  -1.42%  "
$ws.Range("E9").Value = "  -0.19%  "
$ws.Range("E10").Value = "  -1.87%  "
$ws.Range("E12").Value = "  +2.17%  "
$ws.Range("D13").Value = "3.965.27"
$ws.Range("E13").Value = "  -1.45%  "
$ws.Range("E14").Value = "  +2.19%  "
$ws.Range("E15").Value = "  -0.01%  "
$ws.Range("E16").Value = "  -1.26%  "
$ws.Range("D17").Value = "3.387.60"
$ws.Range("E17").Value = "  -1.47%  "
$ws.Range("D18").Value = "60.966.52"
$ws.Range("E18").Value = "  -0.92%  "
$ws.Range("E19").Value = "  -2.53%  "
$ws.Range("E20").Value = "  -3.14%  "
$ws.Range("E21").Value = "  -4.52%  "
$ws.Range("E22").Value = "  -3.19%  "
$ws.Range("E23").Value = "  -1.67%  "
$ws.Range("E24").Value = "  +0.93%  "
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("E26").Value = "  -4.69%  "
$ws.Range("D27").Value = "3.525.55"
$ws.Range("E27").Value = "  -1.39%  "
$ws.Range("E28").Value = "  -0.57%  "
$ws.Range("E29").Value = "  -0.24%  "
$ws.Range("E30").Value = "  -2.85%  "
$ws.Range("E31").Value = "  -3.37%  "
$ws.Range("E32").Value = "  -1.88%  "
$ws.Range("E33").Value = "  -2.23%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("E35").Value = "  -1.89%  "
$ws.Range("E36").Value = "  -0.52%  "
$ws.Range("E37").Value = "  +0.16%  "
$ws.Range("D38").Value = "3.416.15"
$ws.Range("E38").Value = "  -1.40%  "
$ws.Range("E39").Value = "  -2.67%  "
$ws.Range("E40").Value = "  -4.57%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("E41").Value = "  +2.72%  "
$ws.Range("B42").Value = "Hedera"
$ws.Range("C42").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("E42").Value = "  -1.75%  "
$ws.Range("E44").Value = "  -2.65%  "
$ws.Range("E46").Value = "  -1.39%  "
$ws.Range("E47").Value = "  -3.50%  "
$ws.Range("E48").Value = "  -0.61%  "
$ws.Range("D49").Value = "2.475.90"
$ws.Range("E49").Value = "  -4.57%  "
$ws.Range("E50").Value = "  -1.72%  "
$ws.Range("E51").Value = "  -1.53%  "

# --- Numeric-looking price values that must remain stored as Text (matching source) ---
# Force text storage via NumberFormat "@", then restore the default "Normal" style
# so no stray style index is left on the cell (matches original unstyled cells).
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "571.65"
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "142.01"
$c.Style = "Normal"
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "7.53"
$c.Style = "Normal"
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.395"
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "28.18"
$c.Style = "Normal"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "13.81"
$c.Style = "Normal"
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "8.98"
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "383.68"
$c.Style = "Normal"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "74.44"
$c.Style = "Normal"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.179"
$c.Style = "Normal"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.42"
$c.Style = "Normal"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "23.49"
$c.Style = "Normal"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "6.98"
$c.Style = "Normal"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "167.69"
$c.Style = "Normal"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "4.98"
$c.Style = "Normal"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "27.57"
$c.Style = "Normal"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.0774"
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "42.16"
$c.Style = "Normal"
